# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 3
    16 = 1
    17 = 2
    18 = 1
    19 = 2
    20 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
